$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had two header rows (row1 = partial units row using
# "mation"/"pompes)"/"Hiver"/"Eté"/"Année" leftover shared-string text,
# row2 = a second sub-header row with "(m3/s)"/"(MW)"/"(GWh)" units) followed
# by 12 data rows (rows 3-14). The new layout uses a single, full header row
# (row1, columns A-K) followed immediately by the 12 data rows (rows 2-13).

# Step 1: remove the old second header row (row 2). This shifts all the data
# rows up by one (old row 3 -> new row 2, ..., old row 14 -> new row 13) and
# updates the sheet dimension automatically.
$ws.Rows.Item(2).Delete()

# Step 2: re-write row 1 as a full header spanning A1:K1 with the new column
# names. A1:E1 keep the workbook's default cell format (font Arial 10,
# General), while F1:K1 use the same "Arial 9 / General" look as the rest of
# the header/label cells in the sheet.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 must fall back to the default "Normal" look (Arial 10) - E1 carried
# over the old sub-header's "Arial 9" formatting, so reset it explicitly.
$ws.Range("A1:E1").Font.Name = "Arial"
$ws.Range("A1:E1").Font.Size = 10

# F1:K1 use the "Arial 9" label formatting (same visual style used by the
# other text labels in the sheet, e.g. the plant-name column).
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Step 3: update the active selection to match the edited workbook (row 2,
# the first data row, columns A-K).
$ws.Range("A2:K2").Select()
